$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 38695.785
$ws.Range("I28").Value = 51333.2
$ws.Range("J28").Value = 7102.25
$ws.Range("K28").Value = 51333.2
$ws.Range("L28").Value = 7102.25
$ws.Range("M28").Value = -50848.2
$ws.Range("N28").Value = -8072.25
$ws.Range("H98").Value = 1739.3077
$ws.Range("I98").Value = 1798.92
$ws.Range("K98").Value = 1798.92
$ws.Range("M98").Value = -300.9200000000001
$ws.Range("H103").Value = 1055
$ws.Range("I103").Value = 790.1667
$ws.Range("K103").Value = 2370.5001
$ws.Range("M103").Value = -1784.5001
$ws.Range("H111").Value = 30117.818
$ws.Range("I111").Value = 15915.667
$ws.Range("J111").Value = 47160.4
$ws.Range("K111").Value = 47747.001
$ws.Range("L111").Value = 141481.2
$ws.Range("M111").Value = -44680.001
$ws.Range("N111").Value = -147615.2
$ws.Range("H122").Value = 1739.3077
$ws.Range("I122").Value = 1798.92
$ws.Range("K122").Value = 5396.76
$ws.Range("M122").Value = -2946.76
$ws.Range("H138").Value = 5285.7744
$ws.Range("I138").Value = 4050.111
$ws.Range("J138").Value = 5791.273
$ws.Range("K138").Value = 12150.333
$ws.Range("L138").Value = 17373.819
$ws.Range("M138").Value = -7010.332999999999
$ws.Range("N138").Value = -27653.819
$ws.Range("H141").Value = 999.05884
$ws.Range("I141").Value = 1033.6875
$ws.Range("K141").Value = 3101.0625
$ws.Range("M141").Value = 2078.9375

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2042.7
$ws.Range("I2").Value = 2046
$ws.Range("J2").Value = 2013
$ws.Range("K2").Value = 2046
$ws.Range("L2").Value = 2013
$ws.Range("M2").Value = -1933
$ws.Range("N2").Value = -2239
$ws.Range("H32").Value = 403797.7
$ws.Range("I32").Value = 403797.7
$ws.Range("K32").Value = 403797.7
$ws.Range("M32").Value = -403510.7
$ws.Range("H116").Value = 2042.7
$ws.Range("I116").Value = 2046
$ws.Range("J116").Value = 2013
$ws.Range("K116").Value = 2046
$ws.Range("L116").Value = 2013
$ws.Range("M116").Value = 248
$ws.Range("N116").Value = -6601

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2042.7
$ws.Range("I3").Value = 2046
$ws.Range("J3").Value = 2013
$ws.Range("K3").Value = 2046
$ws.Range("L3").Value = 2013
$ws.Range("M3").Value = -1932
$ws.Range("N3").Value = -2241
$ws.Range("H105").Value = 2295.9
$ws.Range("I105").Value = 2137
$ws.Range("J105").Value = 2666.6667
$ws.Range("K105").Value = 2137
$ws.Range("L105").Value = 2666.6667
$ws.Range("M105").Value = -390
$ws.Range("N105").Value = -6160.6667
$ws.Range("H107").Value = 50040228
$ws.Range("I107").Value = 36305.4
$ws.Range("J107").Value = 100044150
$ws.Range("K107").Value = 36305.4
$ws.Range("L107").Value = 100044150
$ws.Range("M107").Value = -34385.4
$ws.Range("N107").Value = -100047990

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 456371.03
$ws.Range("I58").Value = 1078.0769
$ws.Range("J58").Value = 1114016.5
$ws.Range("K58").Value = 1078.0769
$ws.Range("L58").Value = 1114016.5
$ws.Range("M58").Value = -875.0769
$ws.Range("N58").Value = -1114422.5
$ws.Range("H86").Value = 4847.143
$ws.Range("I86").Value = 3923.1428
$ws.Range("K86").Value = 3923.1428
$ws.Range("M86").Value = -2800.1428
$ws.Range("H89").Value = 4847.143
$ws.Range("I89").Value = 3923.1428
$ws.Range("K89").Value = 19615.714
$ws.Range("M89").Value = -13999.714
$ws.Range("H99").Value = 115026.66
$ws.Range("I99").Value = 151466.56
$ws.Range("K99").Value = 151466.56
$ws.Range("M99").Value = -149968.56
$ws.Range("H126").Value = 115026.66
$ws.Range("I126").Value = 151466.56
$ws.Range("K126").Value = 454399.68
$ws.Range("M126").Value = -451929.68
$ws.Range("H132").Value = 1667964.9
$ws.Range("I132").Value = 1001056.5
$ws.Range("K132").Value = 3003169.5
$ws.Range("M132").Value = -3000639.5
$ws.Range("H134").Value = 2633
$ws.Range("I134").Value = 900
$ws.Range("J134").Value = 3499.5
$ws.Range("K134").Value = 2700
$ws.Range("L134").Value = 10498.5
$ws.Range("M134").Value = -165
$ws.Range("N134").Value = -15568.5
$ws.Range("H136").Value = 456371.03
$ws.Range("I136").Value = 1078.0769
$ws.Range("J136").Value = 1114016.5
$ws.Range("K136").Value = 3234.2307
$ws.Range("L136").Value = 3342049.5
$ws.Range("M136").Value = -684.2307000000001
$ws.Range("N136").Value = -3347149.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 969.1111
$ws.Range("I5").Value = 541.1667
$ws.Range("K5").Value = 1623.5001
$ws.Range("M5").Value = -1511.5001
$ws.Range("H8").Value = 4854.7144
$ws.Range("I8").Value = 4854.7144
$ws.Range("K8").Value = 14564.1432
$ws.Range("M8").Value = -14425.1432
$ws.Range("H12").Value = 371.875
$ws.Range("I12").Value = 975
$ws.Range("J12").Value = 170.83333
$ws.Range("K12").Value = 2925
$ws.Range("L12").Value = 512.49999
$ws.Range("M12").Value = -2752
$ws.Range("N12").Value = -858.49999
$ws.Range("H40").Value = 120.73684
$ws.Range("I40").Value = 52.916668
$ws.Range("J40").Value = 237
$ws.Range("K40").Value = 211.666672
$ws.Range("L40").Value = 948
$ws.Range("M40").Value = -142.666672
$ws.Range("N40").Value = -1086
$ws.Range("H68").Value = 3000
$ws.Range("J68").Value = 3000
$ws.Range("L68").Value = 9000
$ws.Range("N68").Value = -10622
$ws.Range("H71").Value = 3000
$ws.Range("J71").Value = 3000
$ws.Range("L71").Value = 27000
$ws.Range("N71").Value = -35112
$ws.Range("H107").Value = 1566.619
$ws.Range("I107").Value = 1246.0667
$ws.Range("J107").Value = 2368
$ws.Range("K107").Value = 3738.2001
$ws.Range("L107").Value = 7104
$ws.Range("M107").Value = -1818.2001
$ws.Range("N107").Value = -10944
$ws.Range("H121").Value = 78293.84
$ws.Range("J121").Value = 126800.375
$ws.Range("L121").Value = 380401.125
$ws.Range("N121").Value = -383021.125
$ws.Range("H122").Value = 335.5
$ws.Range("I122").Value = 261.1111
$ws.Range("J122").Value = 396.36365
$ws.Range("K122").Value = 2349.9999
$ws.Range("L122").Value = 3567.27285
$ws.Range("M122").Value = 100.0000999999997
$ws.Range("N122").Value = -8467.272850000001
$ws.Range("H128").Value = 201749.25
$ws.Range("I128").Value = 201749.25
$ws.Range("K128").Value = 605247.75
$ws.Range("M128").Value = -600267.75
$ws.Range("H131").Value = 11838.823
$ws.Range("I131").Value = 635.8
$ws.Range("J131").Value = 16506.75
$ws.Range("K131").Value = 1907.4
$ws.Range("L131").Value = 49520.25
$ws.Range("M131").Value = 3132.6
$ws.Range("N131").Value = -59600.25
$ws.Range("H135").Value = 969.1111
$ws.Range("I135").Value = 541.1667
$ws.Range("K135").Value = 4870.5003
$ws.Range("M135").Value = -2335.5003

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 10483.823
$ws.Range("I126").Value = 15691.777
$ws.Range("K126").Value = 47075.331
$ws.Range("M126").Value = -44605.331

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4125
$ws.Range("I40").Value = 3250
$ws.Range("K40").Value = 3250
$ws.Range("M40").Value = -3114
$ws.Range("H46").Value = 2594.3572
$ws.Range("I46").Value = 2157.6956
$ws.Range("K46").Value = 2157.6956
$ws.Range("M46").Value = -1969.6956
$ws.Range("H130").Value = 400000
$ws.Range("J130").Value = 400000
$ws.Range("L130").Value = 400000
$ws.Range("N130").Value = -410040
$ws.Range("H136").Value = 5023.702
$ws.Range("I136").Value = 3761.5134
$ws.Range("J136").Value = 9693.799999999999
$ws.Range("K136").Value = 11284.5402
$ws.Range("L136").Value = 29081.4
$ws.Range("M136").Value = -8734.540199999999
$ws.Range("N136").Value = -34181.39999999999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1696.5238
$ws.Range("I126").Value = 1658
$ws.Range("J126").Value = 1927.6666
$ws.Range("K126").Value = 4974
$ws.Range("L126").Value = 5782.9998
$ws.Range("M126").Value = -2504
$ws.Range("N126").Value = -10722.9998
$ws.Range("H132").Value = 265949.62
$ws.Range("I132").Value = 394492.1
$ws.Range("J132").Value = 3722.96
$ws.Range("K132").Value = 1183476.3
$ws.Range("L132").Value = 11168.88
$ws.Range("M132").Value = -1180946.3
$ws.Range("N132").Value = -16228.88
$ws.Range("H136").Value = 2630.425
$ws.Range("I136").Value = 1598.963
$ws.Range("K136").Value = 4796.889
$ws.Range("M136").Value = -2246.889
